$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Totals" header in AG4 -------------------------------------
# Copy formatting from the last year-header cell (AF4) so the new
# header cell matches the existing header styling, then set its text.
$ws.Range("AF4").Copy()
$ws.Range("AG4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("AG4").Value = "Totals"

# "Crop year" -> "Country"
$ws.Range("A4").Value = "Country"

# --- Rename Bolivia ---------------------------------------------------
$ws.Range("A6").Value = "Bolivia"

# --- New "Totals" column (AG) with row sums ---------------------------
# Row 5 gets its own (non-shared) formula.
$ws.Range("AG5").Formula = "=SUM(C5:AF5)"

# Row 61 (the "Total" summary row) must end up using the same
# (non-bold) style as the rest of the column, so strip the bold
# before the fill creates the cell's style.
$ws.Range("AG61").Font.Bold = $false

# Rows 6 through 61 share one formula group (mirrors how Excel
# would fill the formula down the column).
$ws.Range("AG6:AG61").Formula = "=SUM(C6:AF6)"

# Row 60 is a blank separator row - it keeps the new column's number
# format/style but must not carry a formula or value.
$ws.Range("AG60").ClearContents()

# Row 62 (footer / copyright row) also gets the new column's styling
# with no content.
$ws.Range("AG62").NumberFormat = "#,##0"

# --- Selection bookkeeping ---------------------------------------------
$ws.Range("A7").Select()
